# Update "想去人数" (F column) counts and fix a venue name on two sheets:
#   - "展览"   (Exhibitions)
#   - "全部类型" (All types)
# Both sheets list largely the same events, but with slightly different row
# offsets after row 31, so each sheet is updated with its own row numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 9
$ws1.Range("F5").Value  = 13725
$ws1.Range("F7").Value  = 38
$ws1.Range("F10").Value = 129
$ws1.Range("F13").Value = 16
$ws1.Range("F15").Value = 13739
$ws1.Range("F18").Value = 9055
$ws1.Range("F20").Value = 8170
$ws1.Range("F29").Value = 1026
$ws1.Range("F32").Value = 405
$ws1.Range("F33").Value = 5
$ws1.Range("F35").Value = 205

$ws1.Range("D38").Value = "金山南路288号 木渎影视城会展中心"
$ws1.Range("F38").Value = 5032

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 9
$ws4.Range("F5").Value  = 13725
$ws4.Range("F7").Value  = 38
$ws4.Range("F10").Value = 129
$ws4.Range("F13").Value = 16
$ws4.Range("F15").Value = 13739
$ws4.Range("F18").Value = 9055
$ws4.Range("F20").Value = 8170
$ws4.Range("F29").Value = 1026
$ws4.Range("F34").Value = 405
$ws4.Range("F35").Value = 5
$ws4.Range("F37").Value = 205

$ws4.Range("D40").Value = "金山南路288号 木渎影视城会展中心"
$ws4.Range("F40").Value = 5032
